$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 368.7143
$ws.Range("I18").Value = 193.66667
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 193.66667
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = 90.33332999999999
$ws.Range("N18").Value = -1068
$ws.Range("H70").Value = 3882.6924
$ws.Range("H73").Value = 3882.6924
$ws.Range("H82").Value = 9999.625
$ws.Range("I82").Value = 2700
$ws.Range("K82").Value = 8100
$ws.Range("M82").Value = -7694
$ws.Range("H85").Value = 9999.625
$ws.Range("I85").Value = 2700
$ws.Range("K85").Value = 8100
$ws.Range("M85").Value = -6696
$ws.Range("H88").Value = 1834.8462
$ws.Range("I88").Value = 767.6667
$ws.Range("J88").Value = 2155
$ws.Range("K88").Value = 767.6667
$ws.Range("L88").Value = 2155
$ws.Range("M88").Value = -361.6667
$ws.Range("N88").Value = -2967
$ws.Range("H91").Value = 1834.8462
$ws.Range("I91").Value = 767.6667
$ws.Range("J91").Value = 2155
$ws.Range("K91").Value = 767.6667
$ws.Range("L91").Value = 2155
$ws.Range("M91").Value = 636.3333
$ws.Range("N91").Value = -4963
$ws.Range("H129").Value = 867.3
$ws.Range("I129").Value = 397.625
$ws.Range("J129").Value = 908.1413
$ws.Range("K129").Value = 1192.875
$ws.Range("L129").Value = 2724.4239
$ws.Range("M129").Value = 3807.125
$ws.Range("N129").Value = -12724.4239
$ws.Range("H137").Value = 2626.147
$ws.Range("I137").Value = 1473.3334
$ws.Range("J137").Value = 5392.9
$ws.Range("K137").Value = 4420.0002
$ws.Range("L137").Value = 16178.7
$ws.Range("M137").Value = -1870.0002
$ws.Range("N137").Value = -21278.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5317.8096
$ws.Range("I32").Value = 4166
$ws.Range("K32").Value = 4166
$ws.Range("M32").Value = -3879
$ws.Range("H76").Value = 37166.5
$ws.Range("J76").Value = 37166.5
$ws.Range("L76").Value = 37166.5
$ws.Range("N76").Value = -37842.5
$ws.Range("H79").Value = 37166.5
$ws.Range("J79").Value = 37166.5
$ws.Range("L79").Value = 37166.5
$ws.Range("N79").Value = -39506.5
$ws.Range("H132").Value = 1653.4103
$ws.Range("I132").Value = 941.28125
$ws.Range("J132").Value = 4908.857
$ws.Range("K132").Value = 2823.84375
$ws.Range("L132").Value = 14726.571
$ws.Range("M132").Value = -293.84375
$ws.Range("N132").Value = -19786.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 20647
$ws.Range("I6").Value = 20647
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 20647
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -20534
$ws.Range("N6").ClearContents()
$ws.Range("H123").Value = 29920
$ws.Range("J123").Value = 29920
$ws.Range("L123").Value = 29920
$ws.Range("N123").Value = -39720
$ws.Range("H134").Value = 1506.7119
$ws.Range("I134").Value = 1107.2
$ws.Range("K134").Value = 3321.6
$ws.Range("M134").Value = -786.6000000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15154258
$ws.Range("I31").Value = 1355.3
$ws.Range("K31").Value = 1355.3
$ws.Range("M31").Value = -1060.3
$ws.Range("H34").Value = 15154258
$ws.Range("I34").Value = 1355.3
$ws.Range("K34").Value = 1355.3
$ws.Range("M34").Value = -1153.3
$ws.Range("H58").Value = 1571.5238
$ws.Range("I58").Value = 1390.7468
$ws.Range("K58").Value = 1390.7468
$ws.Range("M58").Value = -1187.7468
$ws.Range("H82").Value = 39800
$ws.Range("J82").Value = 39800
$ws.Range("L82").Value = 39800
$ws.Range("N82").Value = -40522
$ws.Range("H85").Value = 39800
$ws.Range("J85").Value = 39800
$ws.Range("L85").Value = 39800
$ws.Range("N85").Value = -42296
$ws.Range("H105").Value = 2558.6667
$ws.Range("I105").Value = 2171.3333
$ws.Range("J105").Value = 3333.3333
$ws.Range("K105").Value = 2171.3333
$ws.Range("L105").Value = 3333.3333
$ws.Range("M105").Value = -424.3332999999998
$ws.Range("N105").Value = -6827.3333
$ws.Range("H134").Value = 6610.7393
$ws.Range("I134").Value = 14973.429
$ws.Range("J134").Value = 2952.0625
$ws.Range("K134").Value = 44920.287
$ws.Range("L134").Value = 8856.1875
$ws.Range("M134").Value = -42385.287
$ws.Range("N134").Value = -13926.1875
$ws.Range("H136").Value = 1571.5238
$ws.Range("I136").Value = 1390.7468
$ws.Range("K136").Value = 4172.2404
$ws.Range("M136").Value = -1622.2404

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2621.7058
$ws.Range("I3").Value = 2332.6
$ws.Range("J3").Value = 4790
$ws.Range("K3").Value = 6997.799999999999
$ws.Range("L3").Value = 14370
$ws.Range("M3").Value = -6885.799999999999
$ws.Range("N3").Value = -14594
$ws.Range("H63").Value = 4937.5
$ws.Range("I63").Value = 4625
$ws.Range("J63").Value = 5250
$ws.Range("K63").Value = 13875
$ws.Range("L63").Value = 15750
$ws.Range("M63").Value = -13126
$ws.Range("N63").Value = -17248
$ws.Range("H66").Value = 4937.5
$ws.Range("I66").Value = 4625
$ws.Range("J66").Value = 5250
$ws.Range("K66").Value = 41625
$ws.Range("L66").Value = 47250
$ws.Range("M66").Value = -37881
$ws.Range("N66").Value = -54738
$ws.Range("H114").Value = 7289.25
$ws.Range("I114").Value = 78.5
$ws.Range("J114").Value = 14500
$ws.Range("K114").Value = 235.5
$ws.Range("L114").Value = 43500
$ws.Range("M114").Value = 3018.5
$ws.Range("N114").Value = -50008
$ws.Range("H117").Value = 8274.5
$ws.Range("I117").Value = 1032.6666
$ws.Range("K117").Value = 3097.9998
$ws.Range("M117").Value = 344.0001999999999
$ws.Range("H121").Value = 2221.9707
$ws.Range("J121").Value = 2315.323
$ws.Range("L121").Value = 6945.968999999999
$ws.Range("N121").Value = -9565.968999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 35988.89
$ws.Range("J15").Value = 35988.89
$ws.Range("L15").Value = 35988.89
$ws.Range("N15").Value = -36564.89
$ws.Range("H64").Value = 29595.4
$ws.Range("J64").Value = 29595.4
$ws.Range("L64").Value = 29595.4
$ws.Range("N64").Value = -30091.4
$ws.Range("H67").Value = 29595.4
$ws.Range("J67").Value = 29595.4
$ws.Range("L67").Value = 29595.4
$ws.Range("N67").Value = -31311.4
$ws.Range("H81").Value = 35988.89
$ws.Range("J81").Value = 35988.89
$ws.Range("L81").Value = 35988.89
$ws.Range("N81").Value = -37984.89
$ws.Range("H84").Value = 35988.89
$ws.Range("J84").Value = 35988.89
$ws.Range("L84").Value = 107966.67
$ws.Range("N84").Value = -117950.67
$ws.Range("H107").Value = 5291346
$ws.Range("I107").Value = 266.6842
$ws.Range("J107").Value = 55556596
$ws.Range("K107").Value = 266.6842
$ws.Range("L107").Value = 55556596
$ws.Range("M107").Value = 1653.3158
$ws.Range("N107").Value = -55560436

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 799.3333
$ws.Range("I68").Value = 799.3333
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 799.3333
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -50.33330000000001
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 799.3333
$ws.Range("I71").Value = 799.3333
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 3996.6665
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -252.6665000000003
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 36000
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("H65").Value = 36000
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("H75").Value = 37250
$ws.Range("J75").Value = 37250
$ws.Range("L75").Value = 37250
$ws.Range("N75").Value = -39122
$ws.Range("H78").Value = 37250
$ws.Range("J78").Value = 37250
$ws.Range("L78").Value = 111750
$ws.Range("N78").Value = -121110
$ws.Range("H126").Value = 2177.9167
$ws.Range("I126").Value = 1532.2142
$ws.Range("J126").Value = 3081.9
$ws.Range("K126").Value = 4596.642599999999
$ws.Range("L126").Value = 9245.700000000001
$ws.Range("M126").Value = -2126.642599999999
$ws.Range("N126").Value = -14185.7
